$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "849÷5=169, 4"
$tbl.Cell(1, 2).Range.Text = "348÷9=38, 6"
$tbl.Cell(1, 3).Range.Text = "945÷5=189, 0"
$tbl.Cell(1, 4).Range.Text = "229÷9=25, 4"
$tbl.Cell(1, 5).Range.Text = "810÷6=135, 0"
$tbl.Cell(5, 1).Range.Text = "566÷6=94, 2"
$tbl.Cell(5, 2).Range.Text = "445÷8=55, 5"
$tbl.Cell(5, 3).Range.Text = "125÷7=17, 6"
$tbl.Cell(5, 4).Range.Text = "563÷4=140, 3"
$tbl.Cell(5, 5).Range.Text = "151÷4=37, 3"
$tbl.Cell(9, 1).Range.Text = "199÷9=22, 1"
$tbl.Cell(9, 2).Range.Text = "226÷9=25, 1"
$tbl.Cell(9, 3).Range.Text = "113÷8=14, 1"
$tbl.Cell(9, 4).Range.Text = "977÷2=488, 1"
$tbl.Cell(9, 5).Range.Text = "392÷9=43, 5"
$tbl.Cell(13, 1).Range.Text = "220÷5=44, 0"
$tbl.Cell(13, 2).Range.Text = "963÷8=120, 3"
$tbl.Cell(13, 3).Range.Text = "741÷8=92, 5"
$tbl.Cell(13, 4).Range.Text = "639÷9=71, 0"
$tbl.Cell(13, 5).Range.Text = "494÷4=123, 2"
$tbl.Cell(17, 1).Range.Text = "776÷7=110, 6"
$tbl.Cell(17, 2).Range.Text = "562÷4=140, 2"
$tbl.Cell(17, 3).Range.Text = "119÷9=13, 2"
$tbl.Cell(17, 4).Range.Text = "768÷3=256, 0"
$tbl.Cell(17, 5).Range.Text = "105÷5=21, 0"
